$p = $ppt.ActivePresentation

# The text run we need to split/update:
#   "C:\Users\user\.tstool\14\DataStores\HydroBase.cfg"
# becomes three separate runs (and fixes the "DataStores" -> "datastores" casing):
#   1) "C:\Users\user"
#   2) "\."
#   3) "tstool\14\datastores\HydroBase.cfg"
$fullTarget = "C:\Users\user\.tstool\14\DataStores\HydroBase.cfg"
$part1 = "C:\Users\user"
$part2 = "\."
$part3 = "tstool\14\datastores\HydroBase.cfg"

# Locate the paragraph (anywhere in the deck) that contains the target text, instead of
# hard-coding slide/shape/paragraph indices.
$targetPara = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $shp = $sl.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text.Contains($fullTarget)) {
                $paraCount = $tr.Paragraphs().Count
                for ($pi = 1; $pi -le $paraCount; $pi++) {
                    $candidate = $tr.Paragraphs($pi, 1)
                    if ($candidate.Text.Contains($fullTarget)) {
                        $targetPara = $candidate
                    }
                }
            }
        }
    }
}

$para = $targetPara

$text = $para.Text
$startOffset = $text.IndexOf($fullTarget)
$start = $startOffset + 1  # TextRange.Characters is 1-based

$start1 = $start
$start2 = $start1 + $part1.Length
$start3 = $start2 + $part2.Length

# Write the last chunk first, then work backwards, so earlier offsets stay valid
# while later ones are still being computed from the original (pre-edit) text.
$sub3 = $para.Characters($start3, $part3.Length)
$sub3.Text = $part3

$sub2 = $para.Characters($start2, $part2.Length)
$sub2.Text = $part2

$sub1 = $para.Characters($start1, $part1.Length)
$sub1.Text = $part1
